# Adds four new weekly columns (24-27) of NRS death data for Health Boards,
# in columns Y:AB, mirroring the existing layout in columns B:X.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): week numbers 24..27 ----
$headerVals = @(24, 25, 26, 27)
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item(1, 25 + $i).Value = $headerVals[$i]
}

# Match header formatting to the existing header cells (e.g. X1), then drop
# the top border that the existing header row carries (the pasted-in columns
# do not have it).
[void]$ws.Range("X1").Copy()
$ws.Range("Y1:AB1").PasteSpecial(-4122)
$ws.Range("Y1:AB1").Borders.Item(3).LineStyle = -4142

# ---- Data rows (2-15): Health Board weekly death counts ----
$data = @(
    @(6, 4, 1, 0),    # row 2  Ayrshire and Arran
    @(3, 1, 0, 1),    # row 3  Borders
    @(0, 0, 0, 0),    # row 4  Dumfries and Galloway
    @(3, 2, 3, 0),    # row 5  Fife
    @(6, 4, 3, 1),    # row 6  Forth Valley
    @(6, 4, 3, 0),    # row 7  Grampian
    @(17, 15, 8, 3),  # row 8  Greater Glasgow and Clyde
    @(0, 0, 0, 0),    # row 9  Highland
    @(8, 8, 11, 2),   # row 10 Lanarkshire
    @(11, 9, 4, 7),   # row 11 Lothian
    @(0, 0, 0, 0),    # row 12 Orkney
    @(0, 0, 0, 0),    # row 13 Shetland
    @(9, 2, 2, 3),    # row 14 Tayside
    @(0, 0, 0, 0)     # row 15 Western Isles
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, 25 + $c).Value = $rowVals[$c]
    }
}

# Formatting for data rows: same Comma-style right aligned numbers as the
# rest of the table (columns B:X).
[void]$ws.Range("X3").Copy()
$ws.Range("Y3:AB15").PasteSpecial(-4122)

[void]$ws.Range("X2").Copy()
$ws.Range("Y2:AB2").PasteSpecial(-4122)

# ---- View state: scroll right so the new columns are visible, and leave
# the selection on the last cell that was filled in. ----
[void]$ws.Range("AB2").Select()
